$d = $word.ActiveDocument

# --- 1. Remove the old "_GoBack" bookmark that currently sits after
#        "Project Management" (Word re-homes this bookmark to the most
#        recent edit location, so it needs to move to the new edit). ---
try {
    $oldGoBack = $d.Bookmarks("_GoBack")
    $oldGoBack.Delete()
} catch {
    # no existing _GoBack bookmark - nothing to remove
}

# --- 2. Locate the "Chase" hobby entry we are about to edit. ---
$rng = $d.Content
$found = $rng.Find.Execute("Chase")

if ($found) {
    # --- 3. Re-create the "_GoBack" bookmark right after this run, i.e.
    #        at the (still-collapsed) end of the found range, before the
    #        text itself is changed. InsertXML on a located range inserts
    #        the supplied fragment at the end of that range without
    #        disturbing the range's own text. ---
    $bookmarkXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p>' +
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
        '<w:bookmarkEnd w:id="0"/>' +
        '</w:p></w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($bookmarkXml)

    # --- 4. Now swap the word itself: Chase -> Chess. ---
    $d.Content.Find.Execute("Chase", $true, $false, $false, $false, $false, `
                             $true, 1, $false, "Chess", 2) | Out-Null
}
